# Update quiz-app subproject status and modify results data structure
# - Replace the sample question/answer rows with a simpler "2+2" row
# - Drop the second (multiple-choice) question row entirely
# - Clear the bold/dark-fill header styling back to the workbook default
# - Remove the explicit (taller) row heights so rows use the sheet default
# - Give columns A and F explicit widths
# - Update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the third data row (the "Capital of France?" question) entirely,
# shifting everything below it up.
$ws.Rows("3:3").Delete()

# Row 2 now holds the only remaining question. Replace its contents.
$ws.Range("A2").Value = "2+2"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3

# Remove the custom header/row styling (bold white-on-dark fill) so the
# cells fall back to the default "Normal" style.
$ws.Range("A1:F2").Style = "Normal"

# Drop the explicit row heights (36 / 24) so rows use the sheet default.
$ws.Rows("1:2").AutoFit()

# Set explicit widths for columns A and F.
$ws.Columns("A").ColumnWidth = 42.36328125
$ws.Columns("F").ColumnWidth = 22

# Update the selected cell shown when the workbook is opened.
$ws.Range("C11").Select()
